$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 200006000
$ws.Range("J69").Value = 200006000
$ws.Range("L69").Value = 600018000
$ws.Range("N69").Value = -600019748
$ws.Range("H72").Value = 200006000
$ws.Range("J72").Value = 200006000
$ws.Range("L72").Value = 1800054000
$ws.Range("N72").Value = -1800062736
$ws.Range("H74").Value = 8030.0713
$ws.Range("I74").Value = 4157.8887
$ws.Range("K74").Value = 4157.8887
$ws.Range("M74").Value = -3221.8887
$ws.Range("H77").Value = 8030.0713
$ws.Range("I77").Value = 4157.8887
$ws.Range("K77").Value = 20789.4435
$ws.Range("M77").Value = -16109.4435
$ws.Range("H138").Value = 5875.378
$ws.Range("J138").Value = 6432.3516
$ws.Range("L138").Value = 19297.0548
$ws.Range("N138").Value = -29577.0548

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3006357
$ws.Range("I8").Value = 5129250
$ws.Range("J8").Value = 175833
$ws.Range("K8").Value = 5129250
$ws.Range("L8").Value = 175833
$ws.Range("M8").Value = -5129106
$ws.Range("N8").Value = -176121
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2350
$ws.Range("H45").Value = 2128
$ws.Range("I45").Value = 1224.1111
$ws.Range("J45").Value = 3483.8333
$ws.Range("K45").Value = 1224.1111
$ws.Range("L45").Value = 3483.8333
$ws.Range("M45").Value = -847.1111000000001
$ws.Range("N45").Value = -4237.8333
$ws.Range("H74").Value = 3456
$ws.Range("I74").Value = 2858.7036
$ws.Range("K74").Value = 2858.7036
$ws.Range("M74").Value = -1984.7036
$ws.Range("H77").Value = 3456
$ws.Range("I77").Value = 2858.7036
$ws.Range("K77").Value = 14293.518
$ws.Range("M77").Value = -9925.518
$ws.Range("H97").Value = 455.33334
$ws.Range("I97").Value = 68.5
$ws.Range("K97").Value = 68.5
$ws.Range("M97").Value = 427.5
$ws.Range("H102").Value = 3351.2173
$ws.Range("I102").Value = 2353.9
$ws.Range("K102").Value = 2353.9
$ws.Range("M102").Value = -731.9000000000001
$ws.Range("H132").Value = 2180
$ws.Range("I132").Value = 1810.5238
$ws.Range("K132").Value = 5431.5714
$ws.Range("M132").Value = -2901.5714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1545
$ws.Range("J5").Value = 612.5
$ws.Range("L5").Value = 612.5
$ws.Range("N5").Value = -838.5
$ws.Range("H81").Value = 21887.5
$ws.Range("J81").Value = 21887.5
$ws.Range("L81").Value = 21887.5
$ws.Range("N81").Value = -24009.5
$ws.Range("H84").Value = 21887.5
$ws.Range("J84").Value = 21887.5
$ws.Range("L84").Value = 65662.5
$ws.Range("N84").Value = -76270.5
$ws.Range("H131").Value = 245713.28
$ws.Range("J131").Value = 241665.67
$ws.Range("L131").Value = 241665.67
$ws.Range("N131").Value = -251745.67

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1738.579
$ws.Range("I132").Value = 1668.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5005.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2475.5
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 2091.5
$ws.Range("I134").Value = 2009.25
$ws.Range("J134").Value = 2749.5
$ws.Range("K134").Value = 6027.75
$ws.Range("L134").Value = 8248.5
$ws.Range("M134").Value = -3492.75
$ws.Range("N134").Value = -13318.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8.666667
$ws.Range("J2").Value = 7.5
$ws.Range("L2").Value = 45
$ws.Range("N2").Value = -271
$ws.Range("H8").Value = 375.25
$ws.Range("I8").Value = 375.25
$ws.Range("K8").Value = 1125.75
$ws.Range("M8").Value = -986.75
$ws.Range("H16").Value = 33.75
$ws.Range("I16").Value = 40
$ws.Range("K16").Value = 120
$ws.Range("M16").Value = 53
$ws.Range("H122").Value = 261.25
$ws.Range("J122").Value = 261.25
$ws.Range("L122").Value = 2351.25
$ws.Range("N122").Value = -7251.25
$ws.Range("H127").Value = 2091.5
$ws.Range("J127").Value = 2091.5
$ws.Range("L127").Value = 6274.5
$ws.Range("N127").Value = -16194.5
$ws.Range("H132").Value = 2314.1428
$ws.Range("I132").Value = 1379.6
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 12416.4
$ws.Range("L132").Value = 25499.9997
$ws.Range("M132").Value = -9886.4
$ws.Range("N132").Value = -30559.9997
$ws.Range("H140").Value = 3335.75
$ws.Range("I140").Value = 3392.111
$ws.Range("K140").Value = 10176.333
$ws.Range("M140").Value = -4996.332999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14028.417
$ws.Range("I80").Value = 20556.834
$ws.Range("J80").Value = 7500
$ws.Range("K80").Value = 20556.834
$ws.Range("L80").Value = 7500
$ws.Range("M80").Value = -19558.834
$ws.Range("N80").Value = -9496
$ws.Range("H83").Value = 14028.417
$ws.Range("I83").Value = 20556.834
$ws.Range("J83").Value = 7500
$ws.Range("K83").Value = 102784.17
$ws.Range("L83").Value = 37500
$ws.Range("M83").Value = -97792.17
$ws.Range("N83").Value = -47484
$ws.Range("H97").Value = 3150.6
$ws.Range("I97").Value = 810.8
$ws.Range("K97").Value = 810.8
$ws.Range("M97").Value = -314.8
$ws.Range("H122").Value = 6404.8965
$ws.Range("I122").Value = 7716.136
$ws.Range("K122").Value = 23148.408
$ws.Range("M122").Value = -20698.408

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3102.3462
$ws.Range("I93").Value = 1427.6
$ws.Range("K93").Value = 1427.6
$ws.Range("M93").Value = -179.5999999999999
$ws.Range("H122").Value = 4130.1113
$ws.Range("I122").Value = 4159
$ws.Range("K122").Value = 12477
$ws.Range("M122").Value = -10027
$ws.Range("H132").Value = 4985.706
$ws.Range("I132").Value = 4685.7036
$ws.Range("K132").Value = 14057.1108
$ws.Range("M132").Value = -11527.1108

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1545.6923
$ws.Range("I100").Value = 1453.909
$ws.Range("K100").Value = 2907.818
$ws.Range("M100").Value = -2366.818
$ws.Range("H107").Value = 2907
$ws.Range("I107").Value = 2907
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 8721
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -6801
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 3975.818
$ws.Range("I122").Value = 4474.222
$ws.Range("K122").Value = 13422.666
$ws.Range("M122").Value = -10972.666
$ws.Range("H126").Value = 3034.1428
$ws.Range("I126").Value = 2959.8462
$ws.Range("K126").Value = 8879.5386
$ws.Range("M126").Value = -6409.5386
$ws.Range("H132").Value = 2145.3125
$ws.Range("I132").Value = 1843.0769
$ws.Range("K132").Value = 5529.2307
$ws.Range("M132").Value = -2999.2307
